# ADD results from server
# Update computed result values in row 2 of each year sheet (2025, 2030, 2035, 2040, 2045, 2050)
# with the latest values returned from the server run.

$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 4480.996542313505
$ws.Range("B2").Value = 1957.662650831837
$ws.Range("E2").Value = 13553.5499643962
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 40342.40040594552
$ws.Range("N2").Value = 5043.751150312321
$ws.Range("O2").Value = 6946.426214796213

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 7441.123814297815
$ws.Range("B2").Value = 12067.56581655979
$ws.Range("E2").Value = 20155.94148885614
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 60721.20459005129
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 6892.460151705974
$ws.Range("O2").Value = 10353.49431872476

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 15683.44064693935
$ws.Range("B2").Value = 13758.73851748809
$ws.Range("E2").Value = 20155.94148885614
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 78590.13056548127
$ws.Range("M2").Value = 3170.265741071549
$ws.Range("N2").Value = 7844.825721754916
$ws.Range("O2").Value = 15725.42758621792

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 15683.44064693935
$ws.Range("B2").Value = 13758.73851748809
$ws.Range("E2").Value = 20155.94148885614
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 78590.13056548127
$ws.Range("M2").Value = 3170.265741071549
$ws.Range("N2").Value = 7844.825721754916
$ws.Range("O2").Value = 15725.42758621792

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 15683.44064693935
$ws.Range("B2").Value = 13758.73851748809
$ws.Range("E2").Value = 20155.94148885614
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 78590.13056548127
$ws.Range("M2").Value = 3170.265741071549
$ws.Range("N2").Value = 7844.825721754916
$ws.Range("O2").Value = 15725.42758621792

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 15683.44064693935
$ws.Range("B2").Value = 13758.73851748809
$ws.Range("E2").Value = 20155.94148885614
$ws.Range("G2").Value = 4231.516049510827
$ws.Range("H2").Value = 48878.76484426508
$ws.Range("I2").Value = 78590.13056548127
$ws.Range("M2").Value = 3170.265741071549
$ws.Range("N2").Value = 7844.825721754916
$ws.Range("O2").Value = 15725.42758621792
